$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 179.125  # H2: 182.33333 -> 179.125
$ws.Cells.Item(2, 9).Value = 179.16667  # I2: 182.33333 -> 179.16667
$ws.Cells.Item(2, 10).Value = 179  # J2: 0 -> 179
$ws.Cells.Item(2, 11).Value = 179.16667  # K2: 182.33333 -> 179.16667
$ws.Cells.Item(2, 12).Value = 179  # L2: 0 -> 179
$ws.Cells.Item(2, 13).Value = -66.16667000000001  # M2: -69.33332999999999 -> -66.16667000000001
$ws.Cells.Item(2, 14).Value = -405  # N2: None -> -405
$ws.Cells.Item(17, 8).Value = 2597.8  # H17: 1665 -> 2597.8
$ws.Cells.Item(17, 10).Value = 2997.25  # J17: 1997.5 -> 2997.25
$ws.Cells.Item(17, 12).Value = 8991.75  # L17: 5992.5 -> 8991.75
$ws.Cells.Item(17, 14).Value = -9327.75  # N17: -6328.5 -> -9327.75
$ws.Cells.Item(33, 8).Value = 471.9  # H33: 509.75 -> 471.9
$ws.Cells.Item(33, 9).Value = 486.77777  # I33: 534.0714 -> 486.77777
$ws.Cells.Item(33, 10).Value = 338  # J33: 339.5 -> 338
$ws.Cells.Item(33, 11).Value = 486.77777  # K33: 534.0714 -> 486.77777
$ws.Cells.Item(33, 12).Value = 338  # L33: 339.5 -> 338
$ws.Cells.Item(33, 13).Value = -257.77777  # M33: -305.0714 -> -257.77777
$ws.Cells.Item(33, 14).Value = -796  # N33: -797.5 -> -796
$ws.Cells.Item(40, 8).Value = 3326.7144  # H40: 3811.4285 -> 3326.7144
$ws.Cells.Item(40, 9).Value = 1633.6666  # I40: 2764.6667 -> 1633.6666
$ws.Cells.Item(40, 11).Value = 1633.6666  # K40: 2764.6667 -> 1633.6666
$ws.Cells.Item(40, 13).Value = -1458.6666  # M40: -2589.6667 -> -1458.6666
$ws.Cells.Item(43, 8).Value = 4794.4  # H43: 4988.4 -> 4794.4
$ws.Cells.Item(43, 9).Value = 4323.3335  # I43: 4646.6665 -> 4323.3335
$ws.Cells.Item(43, 11).Value = 4323.3335  # K43: 4646.6665 -> 4323.3335
$ws.Cells.Item(43, 13).Value = -4254.3335  # M43: -4577.6665 -> -4254.3335
$ws.Cells.Item(69, 8).Value = 16958.334  # H69: 17350 -> 16958.334
$ws.Cells.Item(69, 10).Value = 18750  # J69: 20000 -> 18750
$ws.Cells.Item(69, 12).Value = 56250  # L69: 60000 -> 56250
$ws.Cells.Item(69, 14).Value = -57998  # N69: -61748 -> -57998
$ws.Cells.Item(70, 8).Value = 4756.95  # H70: 5046.8945 -> 4756.95
$ws.Cells.Item(70, 9).Value = 5214.5  # I70: 6237.375 -> 5214.5
$ws.Cells.Item(70, 10).Value = 4299.4  # J70: 4181.091 -> 4299.4
$ws.Cells.Item(70, 11).Value = 15643.5  # K70: 18712.125 -> 15643.5
$ws.Cells.Item(70, 12).Value = 12898.2  # L70: 12543.273 -> 12898.2
$ws.Cells.Item(70, 13).Value = -15373.5  # M70: -18442.125 -> -15373.5
$ws.Cells.Item(70, 14).Value = -13438.2  # N70: -13083.273 -> -13438.2
$ws.Cells.Item(72, 8).Value = 16958.334  # H72: 17350 -> 16958.334
$ws.Cells.Item(72, 10).Value = 18750  # J72: 20000 -> 18750
$ws.Cells.Item(72, 12).Value = 168750  # L72: 180000 -> 168750
$ws.Cells.Item(72, 14).Value = -177486  # N72: -188736 -> -177486
$ws.Cells.Item(73, 8).Value = 4756.95  # H73: 5046.8945 -> 4756.95
$ws.Cells.Item(73, 9).Value = 5214.5  # I73: 6237.375 -> 5214.5
$ws.Cells.Item(73, 10).Value = 4299.4  # J73: 4181.091 -> 4299.4
$ws.Cells.Item(73, 11).Value = 15643.5  # K73: 18712.125 -> 15643.5
$ws.Cells.Item(73, 12).Value = 12898.2  # L73: 12543.273 -> 12898.2
$ws.Cells.Item(73, 13).Value = -14707.5  # M73: -17776.125 -> -14707.5
$ws.Cells.Item(73, 14).Value = -14770.2  # N73: -14415.273 -> -14770.2
$ws.Cells.Item(74, 8).Value = 3933.625  # H74: 3983.625 -> 3933.625
$ws.Cells.Item(74, 9).Value = 3933.625  # I74: 3983.625 -> 3933.625
$ws.Cells.Item(74, 11).Value = 3933.625  # K74: 3983.625 -> 3933.625
$ws.Cells.Item(74, 13).Value = -2997.625  # M74: -3047.625 -> -2997.625
$ws.Cells.Item(77, 8).Value = 3933.625  # H77: 3983.625 -> 3933.625
$ws.Cells.Item(77, 9).Value = 3933.625  # I77: 3983.625 -> 3933.625
$ws.Cells.Item(77, 11).Value = 19668.125  # K77: 19918.125 -> 19668.125
$ws.Cells.Item(77, 13).Value = -14988.125  # M77: -15238.125 -> -14988.125
$ws.Cells.Item(107, 8).Value = 2116.2856  # H107: 1285.8572 -> 2116.2856
$ws.Cells.Item(107, 9).Value = 1878.75  # I107: 989.6667 -> 1878.75
$ws.Cells.Item(107, 10).Value = 2433  # J107: 1508 -> 2433
$ws.Cells.Item(107, 11).Value = 1878.75  # K107: 989.6667 -> 1878.75
$ws.Cells.Item(107, 12).Value = 2433  # L107: 1508 -> 2433
$ws.Cells.Item(107, 13).Value = 41.25  # M107: 930.3333 -> 41.25
$ws.Cells.Item(107, 14).Value = -6273  # N107: -5348 -> -6273
$ws.Cells.Item(113, 8).Value = 1989.6  # H113: 2054.5557 -> 1989.6
$ws.Cells.Item(113, 9).Value = 1989.6  # I113: 2054.5557 -> 1989.6
$ws.Cells.Item(113, 11).Value = 1989.6  # K113: 2054.5557 -> 1989.6
$ws.Cells.Item(113, 13).Value = 1264.4  # M113: 1199.4443 -> 1264.4
$ws.Cells.Item(129, 8).Value = 2267.8667  # H129: 2118.4707 -> 2267.8667
$ws.Cells.Item(129, 9).Value = 1122.2858  # I129: 1094.6666 -> 1122.2858
$ws.Cells.Item(129, 11).Value = 3366.8574  # K129: 3283.9998 -> 3366.8574
$ws.Cells.Item(129, 13).Value = 1633.1426  # M129: 1716.0002 -> 1633.1426
$ws.Cells.Item(132, 8).Value = 5103.5  # H132: 5990.55 -> 5103.5
$ws.Cells.Item(132, 9).Value = 5128.3335  # I132: 6321 -> 5128.3335
$ws.Cells.Item(132, 11).Value = 15385.0005  # K132: 18963 -> 15385.0005
$ws.Cells.Item(132, 13).Value = -12855.0005  # M132: -16433 -> -12855.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4900.1035  # H32: 1674.3334 -> 4900.1035
$ws.Cells.Item(32, 9).Value = 4900.1035  # I32: 1674.3334 -> 4900.1035
$ws.Cells.Item(32, 11).Value = 4900.1035  # K32: 1674.3334 -> 4900.1035
$ws.Cells.Item(32, 13).Value = -4613.1035  # M32: -1387.3334 -> -4613.1035
$ws.Cells.Item(113, 8).Value = 0  # H113: 94500 -> 0
$ws.Cells.Item(113, 10).Value = 0  # J113: 94500 -> 0
$ws.Cells.Item(113, 12).Value = 0  # L113: 94500 -> 0
$ws.Cells.Item(113, 14).ClearContents()  # N113: -103178 -> (removed)
$ws.Cells.Item(132, 8).Value = 5499.75  # H132: 1743.6666 -> 5499.75
$ws.Cells.Item(132, 9).Value = 5499.75  # I132: 1736.6875 -> 5499.75
$ws.Cells.Item(132, 10).Value = 0  # J132: 1799.5 -> 0
$ws.Cells.Item(132, 11).Value = 16499.25  # K132: 5210.0625 -> 16499.25
$ws.Cells.Item(132, 12).Value = 0  # L132: 5398.5 -> 0
$ws.Cells.Item(132, 13).Value = -13969.25  # M132: -2680.0625 -> -13969.25
$ws.Cells.Item(132, 14).ClearContents()  # N132: -10458.5 -> (removed)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(43, 8).Value = 0  # H43: 420000 -> 0
$ws.Cells.Item(43, 10).Value = 0  # J43: 420000 -> 0
$ws.Cells.Item(43, 12).Value = 0  # L43: 420000 -> 0
$ws.Cells.Item(43, 14).ClearContents()  # N43: -420362 -> (removed)
$ws.Cells.Item(107, 8).Value = 2703  # H107: 3000 -> 2703
$ws.Cells.Item(107, 9).Value = 2703  # I107: 3000 -> 2703
$ws.Cells.Item(107, 11).Value = 2703  # K107: 3000 -> 2703
$ws.Cells.Item(107, 13).Value = -783  # M107: -1080 -> -783
$ws.Cells.Item(134, 8).Value = 2833  # H134: 2250 -> 2833
$ws.Cells.Item(134, 9).Value = 2500  # I134: 2250 -> 2500
$ws.Cells.Item(134, 10).Value = 2999.5  # J134: 0 -> 2999.5
$ws.Cells.Item(134, 11).Value = 7500  # K134: 6750 -> 7500
$ws.Cells.Item(134, 12).Value = 8998.5  # L134: 0 -> 8998.5
$ws.Cells.Item(134, 13).Value = -4965  # M134: -4215 -> -4965
$ws.Cells.Item(134, 14).Value = -14068.5  # N134: None -> -14068.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2187.1667  # H16: 2201.6667 -> 2187.1667
$ws.Cells.Item(16, 9).Value = 624.6  # I16: 639.4 -> 624.6
$ws.Cells.Item(16, 10).Value = 10000  # J16: 10013 -> 10000
$ws.Cells.Item(16, 11).Value = 624.6  # K16: 639.4 -> 624.6
$ws.Cells.Item(16, 12).Value = 10000  # L16: 10013 -> 10000
$ws.Cells.Item(16, 13).Value = -337.6  # M16: -352.4 -> -337.6
$ws.Cells.Item(16, 14).Value = -10574  # N16: -10587 -> -10574
$ws.Cells.Item(22, 8).Value = 5001487  # H22: 4446055.5 -> 5001487
$ws.Cells.Item(22, 9).Value = 1624.75  # I22: 1499.7142 -> 1624.75
$ws.Cells.Item(22, 10).Value = 10001349  # J22: 20002000 -> 10001349
$ws.Cells.Item(22, 11).Value = 1624.75  # K22: 1499.7142 -> 1624.75
$ws.Cells.Item(22, 12).Value = 10001349  # L22: 20002000 -> 10001349
$ws.Cells.Item(22, 13).Value = -1274.75  # M22: -1149.7142 -> -1274.75
$ws.Cells.Item(22, 14).Value = -10002049  # N22: -20002700 -> -10002049
$ws.Cells.Item(86, 8).Value = 3959.125  # H86: 3923.1 -> 3959.125
$ws.Cells.Item(86, 9).Value = 3982.1667  # I86: 3953.1428 -> 3982.1667
$ws.Cells.Item(86, 10).Value = 3890  # J86: 3853 -> 3890
$ws.Cells.Item(86, 11).Value = 3982.1667  # K86: 3953.1428 -> 3982.1667
$ws.Cells.Item(86, 12).Value = 3890  # L86: 3853 -> 3890
$ws.Cells.Item(86, 13).Value = -2859.1667  # M86: -2830.1428 -> -2859.1667
$ws.Cells.Item(86, 14).Value = -6136  # N86: -6099 -> -6136
$ws.Cells.Item(89, 8).Value = 3959.125  # H89: 3923.1 -> 3959.125
$ws.Cells.Item(89, 9).Value = 3982.1667  # I89: 3953.1428 -> 3982.1667
$ws.Cells.Item(89, 10).Value = 3890  # J89: 3853 -> 3890
$ws.Cells.Item(89, 11).Value = 19910.8335  # K89: 19765.714 -> 19910.8335
$ws.Cells.Item(89, 12).Value = 19450  # L89: 19265 -> 19450
$ws.Cells.Item(89, 13).Value = -14294.8335  # M89: -14149.714 -> -14294.8335
$ws.Cells.Item(89, 14).Value = -30682  # N89: -30497 -> -30682
$ws.Cells.Item(113, 8).Value = 2187.1667  # H113: 2201.6667 -> 2187.1667
$ws.Cells.Item(113, 9).Value = 624.6  # I113: 639.4 -> 624.6
$ws.Cells.Item(113, 10).Value = 10000  # J113: 10013 -> 10000
$ws.Cells.Item(113, 11).Value = 624.6  # K113: 639.4 -> 624.6
$ws.Cells.Item(113, 12).Value = 10000  # L113: 10013 -> 10000
$ws.Cells.Item(113, 13).Value = 1545.4  # M113: 1530.6 -> 1545.4
$ws.Cells.Item(113, 14).Value = -14340  # N113: -14353 -> -14340
$ws.Cells.Item(134, 8).Value = 2180.8333  # H134: 1946.4286 -> 2180.8333
$ws.Cells.Item(134, 9).Value = 2117  # I134: 1854.1666 -> 2117
$ws.Cells.Item(134, 11).Value = 6351  # K134: 5562.4998 -> 6351
$ws.Cells.Item(134, 13).Value = -3816  # M134: -3027.4998 -> -3816

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 16  # H11: 277.25 -> 16
$ws.Cells.Item(11, 9).Value = 16  # I11: 277.25 -> 16
$ws.Cells.Item(11, 11).Value = 48  # K11: 831.75 -> 48
$ws.Cells.Item(11, 13).Value = 92  # M11: -691.75 -> 92
$ws.Cells.Item(15, 8).Value = 337  # H15: 194.16667 -> 337
$ws.Cells.Item(15, 9).Value = 366.25  # I15: 186.66667 -> 366.25
$ws.Cells.Item(15, 10).Value = 220  # J15: 201.66667 -> 220
$ws.Cells.Item(15, 11).Value = 1098.75  # K15: 560.00001 -> 1098.75
$ws.Cells.Item(15, 12).Value = 660  # L15: 605.00001 -> 660
$ws.Cells.Item(15, 13).Value = -958.75  # M15: -420.00001 -> -958.75
$ws.Cells.Item(15, 14).Value = -940  # N15: -885.00001 -> -940
$ws.Cells.Item(38, 8).Value = 71.59999999999999  # H38: 62.833332 -> 71.59999999999999
$ws.Cells.Item(38, 9).Value = 84.75  # I38: 71.59999999999999 -> 84.75
$ws.Cells.Item(38, 11).Value = 254.25  # K38: 214.8 -> 254.25
$ws.Cells.Item(38, 13).Value = 92.75  # M38: 132.2 -> 92.75
$ws.Cells.Item(49, 8).Value = 2874.75  # H49: 2733 -> 2874.75
$ws.Cells.Item(49, 9).Value = 3333  # I49: 2999.5 -> 3333
$ws.Cells.Item(49, 10).Value = 1500  # J49: 2200 -> 1500
$ws.Cells.Item(49, 11).Value = 9999  # K49: 8998.5 -> 9999
$ws.Cells.Item(49, 12).Value = 4500  # L49: 6600 -> 4500
$ws.Cells.Item(49, 13).Value = -9843  # M49: -8842.5 -> -9843
$ws.Cells.Item(49, 14).Value = -4812  # N49: -6912 -> -4812
$ws.Cells.Item(139, 8).Value = 9258.875  # H139: 9218.556 -> 9258.875
$ws.Cells.Item(139, 9).Value = 9258.875  # I139: 9218.556 -> 9258.875
$ws.Cells.Item(139, 11).Value = 27776.625  # K139: 27655.668 -> 27776.625
$ws.Cells.Item(139, 13).Value = -22636.625  # M139: -22515.668 -> -22636.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 7295551  # H14: 8025106 -> 7295551
$ws.Cells.Item(14, 10).Value = 121.25  # J14: 161.66667 -> 121.25
$ws.Cells.Item(14, 12).Value = 121.25  # L14: 161.66667 -> 121.25
$ws.Cells.Item(14, 14).Value = -457.25  # N14: -497.66667 -> -457.25
$ws.Cells.Item(102, 8).Value = 2499.5  # H102: 1284.4286 -> 2499.5
$ws.Cells.Item(102, 9).Value = 1999  # I102: 1387.5 -> 1999
$ws.Cells.Item(102, 10).Value = 3000  # J102: 666 -> 3000
$ws.Cells.Item(102, 11).Value = 1999  # K102: 1387.5 -> 1999
$ws.Cells.Item(102, 12).Value = 3000  # L102: 666 -> 3000
$ws.Cells.Item(102, 13).Value = -377  # M102: 234.5 -> -377
$ws.Cells.Item(102, 14).Value = -6244  # N102: -3910 -> -6244
$ws.Cells.Item(122, 8).Value = 4326.1113  # H122: 4554.6665 -> 4326.1113
$ws.Cells.Item(122, 9).Value = 3387.4  # I122: 3798.8 -> 3387.4
$ws.Cells.Item(122, 11).Value = 10162.2  # K122: 11396.4 -> 10162.2
$ws.Cells.Item(122, 13).Value = -7712.200000000001  # M122: -8946.400000000001 -> -7712.200000000001
$ws.Cells.Item(132, 8).Value = 3263.75  # H132: 2412.2104 -> 3263.75
$ws.Cells.Item(132, 9).Value = 3241  # I132: 2411.3333 -> 3241
$ws.Cells.Item(132, 10).Value = 3332  # J132: 2415.5 -> 3332
$ws.Cells.Item(132, 11).Value = 9723  # K132: 7233.999899999999 -> 9723
$ws.Cells.Item(132, 12).Value = 9996  # L132: 7246.5 -> 9996
$ws.Cells.Item(132, 13).Value = -7193  # M132: -4703.999899999999 -> -7193
$ws.Cells.Item(132, 14).Value = -15056  # N132: -12306.5 -> -15056

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 4610.364  # H22: 4993 -> 4610.364
$ws.Cells.Item(22, 9).Value = 3165.6  # I22: 3189.4 -> 3165.6
$ws.Cells.Item(22, 10).Value = 5814.3335  # J22: 7999 -> 5814.3335
$ws.Cells.Item(22, 11).Value = 3165.6  # K22: 3189.4 -> 3165.6
$ws.Cells.Item(22, 12).Value = 5814.3335  # L22: 7999 -> 5814.3335
$ws.Cells.Item(22, 13).Value = -2870.6  # M22: -2894.4 -> -2870.6
$ws.Cells.Item(22, 14).Value = -6404.3335  # N22: -8589 -> -6404.3335
$ws.Cells.Item(27, 8).Value = 4610.364  # H27: 4993 -> 4610.364
$ws.Cells.Item(27, 9).Value = 3165.6  # I27: 3189.4 -> 3165.6
$ws.Cells.Item(27, 10).Value = 5814.3335  # J27: 7999 -> 5814.3335
$ws.Cells.Item(27, 11).Value = 3165.6  # K27: 3189.4 -> 3165.6
$ws.Cells.Item(27, 12).Value = 5814.3335  # L27: 7999 -> 5814.3335
$ws.Cells.Item(27, 13).Value = -3058.6  # M27: -3082.4 -> -3058.6
$ws.Cells.Item(27, 14).Value = -6028.3335  # N27: -8213 -> -6028.3335
$ws.Cells.Item(35, 8).Value = 1556.1428  # H35: 1575 -> 1556.1428
$ws.Cells.Item(35, 9).Value = 1415.5  # I35: 1300 -> 1415.5
$ws.Cells.Item(35, 11).Value = 1415.5  # K35: 1300 -> 1415.5
$ws.Cells.Item(35, 13).Value = -1079.5  # M35: -964 -> -1079.5
$ws.Cells.Item(61, 8).Value = 2707.875  # H61: 2858.3333 -> 2707.875
$ws.Cells.Item(61, 9).Value = 2611.6667  # I61: 2789.25 -> 2611.6667
$ws.Cells.Item(61, 11).Value = 2611.6667  # K61: 2789.25 -> 2611.6667
$ws.Cells.Item(61, 13).Value = -2409.6667  # M61: -2587.25 -> -2409.6667
$ws.Cells.Item(98, 8).Value = 51128.332  # H98: 52015.5 -> 51128.332
$ws.Cells.Item(98, 10).Value = 51128.332  # J98: 52015.5 -> 51128.332
$ws.Cells.Item(98, 12).Value = 51128.332  # L98: 52015.5 -> 51128.332
$ws.Cells.Item(98, 14).Value = -57118.332  # N98: -58005.5 -> -57118.332
$ws.Cells.Item(113, 8).Value = 2707.875  # H113: 2858.3333 -> 2707.875
$ws.Cells.Item(113, 9).Value = 2611.6667  # I113: 2789.25 -> 2611.6667
$ws.Cells.Item(113, 11).Value = 2611.6667  # K113: 2789.25 -> 2611.6667
$ws.Cells.Item(113, 13).Value = -441.6667000000002  # M113: -619.25 -> -441.6667000000002
$ws.Cells.Item(136, 8).Value = 20083.334  # H136: 20092.777 -> 20083.334
$ws.Cells.Item(136, 9).Value = 18764.8  # I136: 18781.8 -> 18764.8
$ws.Cells.Item(136, 11).Value = 56294.39999999999  # K136: 56345.39999999999 -> 56294.39999999999
$ws.Cells.Item(136, 13).Value = -53744.39999999999  # M136: -53795.39999999999 -> -53744.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 14437000  # H2: 16844000 -> 14437000
$ws.Cells.Item(2, 9).Value = 20010800  # I2: 25011000 -> 20010800
$ws.Cells.Item(2, 10).Value = 502500  # J2: 510000.5 -> 502500
$ws.Cells.Item(2, 11).Value = 20010800  # K2: 25011000 -> 20010800
$ws.Cells.Item(2, 12).Value = 502500  # L2: 510000.5 -> 502500
$ws.Cells.Item(2, 13).Value = -20010688  # M2: -25010888 -> -20010688
$ws.Cells.Item(2, 14).Value = -502724  # N2: -510224.5 -> -502724
$ws.Cells.Item(4, 8).Value = 3340476.8  # H4: 3770875 -> 3340476.8
$ws.Cells.Item(4, 9).Value = 1001  # I4: 0 -> 1001
$ws.Cells.Item(4, 10).Value = 3757911.2  # J4: 3770875 -> 3757911.2
$ws.Cells.Item(4, 11).Value = 1001  # K4: 0 -> 1001
$ws.Cells.Item(4, 12).Value = 3757911.2  # L4: 3770875 -> 3757911.2
$ws.Cells.Item(4, 13).Value = -888  # M4: None -> -888
$ws.Cells.Item(4, 14).Value = -3758137.2  # N4: -3771101 -> -3758137.2
$ws.Cells.Item(97, 8).Value = 26664.666  # H97: 29998 -> 26664.666
$ws.Cells.Item(97, 10).Value = 26664.666  # J97: 29998 -> 26664.666
$ws.Cells.Item(97, 12).Value = 26664.666  # L97: 29998 -> 26664.666
$ws.Cells.Item(97, 14).Value = -28646.666  # N97: -31980 -> -28646.666
$ws.Cells.Item(98, 8).Value = 38000  # H98: 38295 -> 38000
$ws.Cells.Item(98, 10).Value = 36000  # J98: 36590 -> 36000
$ws.Cells.Item(98, 12).Value = 36000  # L98: 36590 -> 36000
$ws.Cells.Item(98, 14).Value = -41990  # N98: -42580 -> -41990
$ws.Cells.Item(100, 8).Value = 7165.222  # H100: 7165.3335 -> 7165.222
$ws.Cells.Item(100, 9).Value = 5212.5713  # I100: 5212.7144 -> 5212.5713
$ws.Cells.Item(100, 11).Value = 10425.1426  # K100: 10425.4288 -> 10425.1426
$ws.Cells.Item(100, 13).Value = -9884.142599999999  # M100: -9884.4288 -> -9884.142599999999
$ws.Cells.Item(122, 8).Value = 3539.125  # H122: 2596.9285 -> 3539.125
$ws.Cells.Item(122, 9).Value = 3329.7144  # I122: 2411.6924 -> 3329.7144
$ws.Cells.Item(122, 11).Value = 9989.143199999999  # K122: 7235.0772 -> 9989.143199999999
$ws.Cells.Item(122, 13).Value = -7539.143199999999  # M122: -4785.0772 -> -7539.143199999999
$ws.Cells.Item(136, 8).Value = 13539.429  # H136: 14694.385 -> 13539.429
$ws.Cells.Item(136, 9).Value = 10471.083  # I136: 11557.091 -> 10471.083
$ws.Cells.Item(136, 11).Value = 31413.249  # K136: 34671.273 -> 31413.249
$ws.Cells.Item(136, 13).Value = -28863.249  # M136: -32121.273 -> -28863.249
